# Auto-generated Excel COM-interop script
# Applies scheduled market-data / profit-calc updates to the Kujata profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Cells.Item(4, 8).Value = 301
$ws.Cells.Item(4, 9).Value = 301
$ws.Cells.Item(4, 11).Value = 301
$ws.Cells.Item(4, 13).Value = -187
# Row 40
$ws.Cells.Item(40, 8).Value = 3533.6667
$ws.Cells.Item(40, 9).Value = 3459.8
$ws.Cells.Item(40, 10).Value = 3626
$ws.Cells.Item(40, 11).Value = 3459.8
$ws.Cells.Item(40, 12).Value = 3626
$ws.Cells.Item(40, 13).Value = -3284.8
$ws.Cells.Item(40, 14).Value = -3976
# Row 98
$ws.Cells.Item(98, 8).Value = 1422.579
$ws.Cells.Item(98, 9).Value = 1401.6111
$ws.Cells.Item(98, 10).Value = 1800
$ws.Cells.Item(98, 11).Value = 1401.6111
$ws.Cells.Item(98, 12).Value = 1800
$ws.Cells.Item(98, 13).Value = 96.38889999999992
$ws.Cells.Item(98, 14).Value = -4796
# Row 122
$ws.Cells.Item(122, 8).Value = 1422.579
$ws.Cells.Item(122, 9).Value = 1401.6111
$ws.Cells.Item(122, 10).Value = 1800
$ws.Cells.Item(122, 11).Value = 4204.8333
$ws.Cells.Item(122, 12).Value = 5400
$ws.Cells.Item(122, 13).Value = -1754.8333
$ws.Cells.Item(122, 14).Value = -10300
# Row 137
$ws.Cells.Item(137, 8).Value = 1376.9615
$ws.Cells.Item(137, 9).Value = 875
$ws.Cells.Item(137, 10).Value = 2180.1
$ws.Cells.Item(137, 11).Value = 2625
$ws.Cells.Item(137, 12).Value = 6540.299999999999
$ws.Cells.Item(137, 13).Value = -75
$ws.Cells.Item(137, 14).Value = -11640.3

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 3017.4443
$ws.Cells.Item(32, 9).Value = 3027.36
$ws.Cells.Item(32, 11).Value = 3027.36
$ws.Cells.Item(32, 13).Value = -2740.36
# Row 61
$ws.Cells.Item(61, 8).Value = 1255.8334
$ws.Cells.Item(61, 9).Value = 962.4167
$ws.Cells.Item(61, 10).Value = 1842.6666
$ws.Cells.Item(61, 11).Value = 962.4167
$ws.Cells.Item(61, 12).Value = 1842.6666
$ws.Cells.Item(61, 13).Value = -750.4167
$ws.Cells.Item(61, 14).Value = -2266.6666
# Row 74
$ws.Cells.Item(74, 8).Value = 1116.5714
$ws.Cells.Item(74, 9).Value = 974.6667
$ws.Cells.Item(74, 11).Value = 974.6667
$ws.Cells.Item(74, 13).Value = -100.6667
# Row 77
$ws.Cells.Item(77, 8).Value = 1116.5714
$ws.Cells.Item(77, 9).Value = 974.6667
$ws.Cells.Item(77, 11).Value = 4873.3335
$ws.Cells.Item(77, 13).Value = -505.3334999999997
# Row 132
$ws.Cells.Item(132, 8).Value = 3668.261
$ws.Cells.Item(132, 9).Value = 3445.4119
$ws.Cells.Item(132, 10).Value = 4299.6665
$ws.Cells.Item(132, 11).Value = 10336.2357
$ws.Cells.Item(132, 12).Value = 12898.9995
$ws.Cells.Item(132, 13).Value = -7806.235700000001
$ws.Cells.Item(132, 14).Value = -17958.9995
# Row 136
$ws.Cells.Item(136, 8).Value = 1255.8334
$ws.Cells.Item(136, 9).Value = 962.4167
$ws.Cells.Item(136, 10).Value = 1842.6666
$ws.Cells.Item(136, 11).Value = 2887.2501
$ws.Cells.Item(136, 12).Value = 5527.9998
$ws.Cells.Item(136, 13).Value = -337.2501000000002
$ws.Cells.Item(136, 14).Value = -10627.9998
# Row 37
$ws.Cells.Item(37, 8).Value = 28000
$ws.Cells.Item(37, 10).Value = 28000
$ws.Cells.Item(37, 12).Value = 28000
$ws.Cells.Item(37, 14).Value = -28546

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 1565.3462
$ws.Cells.Item(20, 9).Value = 1299
$ws.Cells.Item(20, 11).Value = 1299
$ws.Cells.Item(20, 13).Value = -1052
# Row 24
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).ClearContents()
$ws.Cells.Item(24, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 1360.375
$ws.Cells.Item(31, 9).Value = 828.05884
$ws.Cells.Item(31, 10).Value = 1963.6666
$ws.Cells.Item(31, 11).Value = 828.05884
$ws.Cells.Item(31, 12).Value = 1963.6666
$ws.Cells.Item(31, 13).Value = -533.05884
$ws.Cells.Item(31, 14).Value = -2553.6666
# Row 34
$ws.Cells.Item(34, 8).Value = 1360.375
$ws.Cells.Item(34, 9).Value = 828.05884
$ws.Cells.Item(34, 10).Value = 1963.6666
$ws.Cells.Item(34, 11).Value = 828.05884
$ws.Cells.Item(34, 12).Value = 1963.6666
$ws.Cells.Item(34, 13).Value = -626.05884
$ws.Cells.Item(34, 14).Value = -2367.6666
# Row 132
$ws.Cells.Item(132, 8).Value = 4771.5
$ws.Cells.Item(132, 9).Value = 5241.6206
$ws.Cells.Item(132, 11).Value = 15724.8618
$ws.Cells.Item(132, 13).Value = -13194.8618
# Row 7
$ws.Cells.Item(7, 8).Value = 259.8
$ws.Cells.Item(7, 9).Value = 259.8
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 259.8
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = -146.8
$ws.Cells.Item(7, 14).ClearContents()
# Row 103
$ws.Cells.Item(103, 8).Value = 2350
$ws.Cells.Item(103, 9).Value = 2350
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 11).Value = 2350
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 13).Value = -1178
$ws.Cells.Item(103, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Cells.Item(68, 8).Value = 1738.9642
$ws.Cells.Item(68, 9).Value = 743.1875
$ws.Cells.Item(68, 10).Value = 3066.6667
$ws.Cells.Item(68, 11).Value = 2229.5625
$ws.Cells.Item(68, 12).Value = 9200.000100000001
$ws.Cells.Item(68, 13).Value = -1418.5625
$ws.Cells.Item(68, 14).Value = -10822.0001
# Row 71
$ws.Cells.Item(71, 8).Value = 1738.9642
$ws.Cells.Item(71, 9).Value = 743.1875
$ws.Cells.Item(71, 10).Value = 3066.6667
$ws.Cells.Item(71, 11).Value = 6688.6875
$ws.Cells.Item(71, 12).Value = 27600.0003
$ws.Cells.Item(71, 13).Value = -2632.6875
$ws.Cells.Item(71, 14).Value = -35712.0003
# Row 131
$ws.Cells.Item(131, 8).Value = 23811038
$ws.Cells.Item(131, 10).Value = 1627.8611
$ws.Cells.Item(131, 12).Value = 4883.5833
$ws.Cells.Item(131, 14).Value = -14963.5833

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 40914320
$ws.Cells.Item(70, 9).Value = 41671570
$ws.Cells.Item(70, 10).Value = 40005620
$ws.Cells.Item(70, 11).Value = 41671570
$ws.Cells.Item(70, 12).Value = 40005620
$ws.Cells.Item(70, 13).Value = -41671300
$ws.Cells.Item(70, 14).Value = -40006160
# Row 73
$ws.Cells.Item(73, 8).Value = 40914320
$ws.Cells.Item(73, 9).Value = 41671570
$ws.Cells.Item(73, 10).Value = 40005620
$ws.Cells.Item(73, 11).Value = 41671570
$ws.Cells.Item(73, 12).Value = 40005620
$ws.Cells.Item(73, 13).Value = -41670634
$ws.Cells.Item(73, 14).Value = -40007492
# Row 80
$ws.Cells.Item(80, 8).Value = 2733.611
$ws.Cells.Item(80, 9).Value = 1527.2727
$ws.Cells.Item(80, 10).Value = 4629.2856
$ws.Cells.Item(80, 11).Value = 1527.2727
$ws.Cells.Item(80, 12).Value = 4629.2856
$ws.Cells.Item(80, 13).Value = -529.2727
$ws.Cells.Item(80, 14).Value = -6625.2856
# Row 83
$ws.Cells.Item(83, 8).Value = 2733.611
$ws.Cells.Item(83, 9).Value = 1527.2727
$ws.Cells.Item(83, 10).Value = 4629.2856
$ws.Cells.Item(83, 11).Value = 7636.363499999999
$ws.Cells.Item(83, 12).Value = 23146.428
$ws.Cells.Item(83, 13).Value = -2644.363499999999
$ws.Cells.Item(83, 14).Value = -33130.428
# Row 97
$ws.Cells.Item(97, 8).Value = 1023
$ws.Cells.Item(97, 9).Value = 1151.5
$ws.Cells.Item(97, 10).Value = 830.25
$ws.Cells.Item(97, 11).Value = 1151.5
$ws.Cells.Item(97, 12).Value = 830.25
$ws.Cells.Item(97, 13).Value = -655.5
$ws.Cells.Item(97, 14).Value = -1822.25
# Row 122
$ws.Cells.Item(122, 8).Value = 3548.3845
$ws.Cells.Item(122, 9).Value = 2974.4285
$ws.Cells.Item(122, 10).Value = 4218
$ws.Cells.Item(122, 11).Value = 8923.2855
$ws.Cells.Item(122, 12).Value = 12654
$ws.Cells.Item(122, 13).Value = -6473.2855
$ws.Cells.Item(122, 14).Value = -17554
# Row 132
$ws.Cells.Item(132, 8).Value = 3055.9443
$ws.Cells.Item(132, 9).Value = 2705.9
$ws.Cells.Item(132, 10).Value = 3493.5
$ws.Cells.Item(132, 11).Value = 8117.700000000001
$ws.Cells.Item(132, 12).Value = 10480.5
$ws.Cells.Item(132, 13).Value = -5587.700000000001
$ws.Cells.Item(132, 14).Value = -15540.5
# Row 135
$ws.Cells.Item(135, 8).Value = 49997.5
$ws.Cells.Item(135, 10).Value = 49997.5
$ws.Cells.Item(135, 12).Value = 49997.5
$ws.Cells.Item(135, 14).Value = -60137.5

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Cells.Item(2, 8).Value = 2805.818
$ws.Cells.Item(2, 10).Value = 2888.4
$ws.Cells.Item(2, 12).Value = 2888.4
$ws.Cells.Item(2, 14).Value = -3112.4
# Row 7
$ws.Cells.Item(7, 8).Value = 1973.2858
$ws.Cells.Item(7, 9).Value = 1561.6
$ws.Cells.Item(7, 11).Value = 1561.6
$ws.Cells.Item(7, 13).Value = -1449.6
# Row 22
$ws.Cells.Item(22, 8).Value = 1999.5
$ws.Cells.Item(22, 10).Value = 1999.5
$ws.Cells.Item(22, 12).Value = 1999.5
$ws.Cells.Item(22, 14).Value = -2589.5
# Row 27
$ws.Cells.Item(27, 8).Value = 1999.5
$ws.Cells.Item(27, 10).Value = 1999.5
$ws.Cells.Item(27, 12).Value = 1999.5
$ws.Cells.Item(27, 14).Value = -2213.5
# Row 46
$ws.Cells.Item(46, 8).Value = 1347.75
$ws.Cells.Item(46, 9).Value = 650
$ws.Cells.Item(46, 10).Value = 1580.3334
$ws.Cells.Item(46, 11).Value = 650
$ws.Cells.Item(46, 12).Value = 1580.3334
$ws.Cells.Item(46, 13).Value = -462
$ws.Cells.Item(46, 14).Value = -1956.3334
# Row 55
$ws.Cells.Item(55, 8).Value = 510.93332
$ws.Cells.Item(55, 9).Value = 73.42856999999999
$ws.Cells.Item(55, 10).Value = 893.75
$ws.Cells.Item(55, 11).Value = 73.42856999999999
$ws.Cells.Item(55, 12).Value = 893.75
$ws.Cells.Item(55, 13).Value = 99.57143000000001
$ws.Cells.Item(55, 14).Value = -1239.75
# Row 68
$ws.Cells.Item(68, 8).Value = 1482
$ws.Cells.Item(68, 9).Value = 1157.7142
$ws.Cells.Item(68, 11).Value = 1157.7142
$ws.Cells.Item(68, 13).Value = -408.7141999999999
# Row 71
$ws.Cells.Item(71, 8).Value = 1482
$ws.Cells.Item(71, 9).Value = 1157.7142
$ws.Cells.Item(71, 11).Value = 5788.571
$ws.Cells.Item(71, 13).Value = -2044.571
# Row 100
$ws.Cells.Item(100, 8).Value = 2125.75
$ws.Cells.Item(100, 9).Value = 2500
$ws.Cells.Item(100, 10).Value = 2001
$ws.Cells.Item(100, 11).Value = 2500
$ws.Cells.Item(100, 12).Value = 2001
$ws.Cells.Item(100, 13).Value = -1959
$ws.Cells.Item(100, 14).Value = -3083
# Row 126
$ws.Cells.Item(126, 8).Value = 1973.2858
$ws.Cells.Item(126, 9).Value = 1561.6
$ws.Cells.Item(126, 11).Value = 4684.799999999999
$ws.Cells.Item(126, 13).Value = -2214.799999999999
# Row 136
$ws.Cells.Item(136, 8).Value = 1991.0625
$ws.Cells.Item(136, 9).Value = 2138.2222
$ws.Cells.Item(136, 11).Value = 6414.6666
$ws.Cells.Item(136, 13).Value = -3864.6666

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Cells.Item(2, 8).Value = 5000
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 14).ClearContents()
